$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Software"
$ws.Range("B10").Value = "imageJ"
$ws.Range("C10").Value = "1.54p"
$ws.Range("D10").Value = "Fluorescence overlay"

$ws.Range("D10").Select()
